# ---------------------------------------------------------------------------
# Atualiza a base "Pediatrico" que alimenta o Power BI:
#   - vira a competencia de 20/10 a 19 / 2026  ->  20/11 a 19/12 / 2025
#   - Nilopolis  passa a ser Itaguai
#   - Paracambi  passa a ser Japeri
#   - Queimados  passa a ser Mage
#   - Sao Joao de Meriti passa a ser Mesquita
#   - inclui novos lancamentos para Queimados e Seropedica
# (tambem serve de "backup" dos dados antigos antes de sobrescrever o
#  arquivo que alimenta o pbi)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ano = "2025"
$competencia = "20/11 a 19/12"

# A coluna "Ano" guarda o ano como texto (ex.: "2025"); forca formato de
# texto para o intervalo que sera escrito antes de preencher, evitando que
# o Excel converta o valor automaticamente para numero.
$ws.Range("F2:F91").NumberFormat = "@"

# Procedimentos cirurgicos (nome + valor unitario), na ordem das colunas B/C
$cirurgias = @(
    @("ADENOIDECTOMIA PEDIÁTRICO", 5330),
    @("AMIGDALECTOMIA- PEDIATRICO", 6713.01),
    @("AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO", 7698.35),
    @("TRATAMENTO CIRÚRGICO DE PERFURAÇÃO DO SEPTO NASAL - PEDIATRICO", 6500),
    @("CORREÇÃO CIRÚRGICO DE ESTRABISMO (ACIMA DE 2 MUSCULOS) - PEDIATRICO", 5255.28),
    @("HERNIOPLASTIA INGUINAL (BILATERAL) - PEDIATRICO", 5850),
    @("HERNIOPLASTIA UMBILICAL - PEDIATRICO", 5237.06),
    @("ORQUIDOPEXIA BILATERAL - PEDIATRICO", 7157.78),
    @("TRATAMENTO CIRÚRGICO DE HIDROCELE - PEDIATRICO", 3782.7),
    @("CORRECAO DE HIPOSPADIA (1º TEMPO) - PEDIATRICO", 6608.86),
    @("PLASTICA TOTAL DO PENIS - PEDIATRICO", 6500),
    @("POSTECTOMIA - PEDIATRICO", 4850)
)

# Consultas (nome + valor unitario), na ordem das colunas I/J
$consultas = @(
    @("CONSULTA PEDIATRICA OTORRINO", 300),
    @("CONSULTA PEDIATRICA CIRURGIA GERAL", 300),
    @("CONSULTA PEDIATRICA OFTALMOLOGISTA", 300)
)

# Cada municipio: nome, quantidades das 12 cirurgias, quantidades das 3 consultas
$municipios = @(
    @("Itaguaí",    @(1, 1, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(5, 0, 0)),
    @("Japeri",     @(1, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(2, 0, 0)),
    @("Magé",       @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(1, 0, 0)),
    @("Mesquita",   @(0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(2, 0, 0)),
    @("Queimados",  @(1, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(3, 0, 0)),
    @("Seropédica", @(0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0), @(1, 0, 0))
)

$linha = 2

for ($m = 0; $m -lt $municipios.Length; $m++) {
    $municipio = $municipios[$m][0]
    $qtdCirurgias = $municipios[$m][1]
    $qtdConsultas = $municipios[$m][2]

    # 12 linhas de cirurgia
    for ($i = 0; $i -lt $cirurgias.Length; $i++) {
        $nomeCirurgia = $cirurgias[$i][0]
        $valorUnitario = $cirurgias[$i][1]
        $qtd = $qtdCirurgias[$i]
        $total = $valorUnitario * $qtd

        $ws.Cells.Item($linha, 1).Value = "Pediatrico"
        $ws.Cells.Item($linha, 2).Value = $nomeCirurgia
        $ws.Cells.Item($linha, 3).Value = $valorUnitario
        $ws.Cells.Item($linha, 4).Value = $qtd
        $ws.Cells.Item($linha, 5).Value = $municipio
        $ws.Cells.Item($linha, 6).Value = $ano
        $ws.Cells.Item($linha, 7).Value = $competencia
        $ws.Cells.Item($linha, 8).Value = $total
        $ws.Cells.Item($linha, 9).Value = ""
        $ws.Cells.Item($linha, 10).Value = ""
        $ws.Cells.Item($linha, 11).Value = ""
        $ws.Cells.Item($linha, 12).Value = ""

        $linha = $linha + 1
    }

    # 3 linhas de consulta
    for ($j = 0; $j -lt $consultas.Length; $j++) {
        $nomeConsulta = $consultas[$j][0]
        $valorConsulta = $consultas[$j][1]
        $qtdC = $qtdConsultas[$j]
        $totalC = $valorConsulta * $qtdC

        $ws.Cells.Item($linha, 1).Value = "Pediatrico"
        $ws.Cells.Item($linha, 2).Value = ""
        $ws.Cells.Item($linha, 3).Value = ""
        $ws.Cells.Item($linha, 4).Value = ""
        $ws.Cells.Item($linha, 5).Value = $municipio
        $ws.Cells.Item($linha, 6).Value = $ano
        $ws.Cells.Item($linha, 7).Value = $competencia
        $ws.Cells.Item($linha, 8).Value = ""
        $ws.Cells.Item($linha, 9).Value = $nomeConsulta
        $ws.Cells.Item($linha, 10).Value = $valorConsulta
        $ws.Cells.Item($linha, 11).Value = $qtdC
        $ws.Cells.Item($linha, 12).Value = $totalC

        $linha = $linha + 1
    }
}

Write-Host "Linhas escritas ate:" ($linha - 1)
